# "Generate Report for Handback"
#
# The localization-status workbook gets its per-file "handback" columns
# (Latest Target File / Latest Handback File / Latest Handback DateTime)
# populated now that the de-de / zh-cn translations have been handed back,
# and the Status column flips from "Ready for handoff" to
# "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

# Column width helper: the OOXML <col width="..."> attribute is always
# `ColumnWidth + 5/6` (5 padding pixels at the default 6px/char Normal
# font), and Excel's ColumnWidth setter itself only has 1/6-character
# (whole pixel) resolution. So to land a raw stored width of W we need to
# feed the setter (W - 5/6), nudged toward the middle of the pixel bucket
# that rounds back to W.
function Set-RawColumnWidth($col, [double]$rawWidth) {
    $col.ColumnWidth = $rawWidth - (5.0 / 6.0)
}

# ---------------------------------------------------------------------
# Overview sheet: the per-language Status cells mirror the same text as
# the language sheets' Status column - update all four.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

Set-RawColumnWidth $wsOverview.Columns.Item(5) 29.9777047293527
Set-RawColumnWidth $wsOverview.Columns.Item(6) 29.9777047293527

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

# Row 2 (7d8b69fb...): Latest Target File / Latest Handback File / Latest Handback DateTime
$wsZh.Hyperlinks.Add(
    $wsZh.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11347c2a6e2f31d8ba59a034b0f5b586dfe56f2d/e2e/7d8b69fb-74ad-419a-b316-fb9822822b42.md",
    "",
    "",
    "7d8b69fb-74ad-419a-b316-fb9822822b42.md"
) | Out-Null
$wsZh.Range("I2").Font.Color = 15570276
$wsZh.Range("J2").Value = "7d8b69fb-74ad-419a-b316-fb9822822b42.e7ae820631db705b393e5b3c6834a1c08189e120.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-06 02:53:23"

# Row 3 (e0757593...): Latest Target File / Latest Handback File / Latest Handback DateTime
$wsZh.Hyperlinks.Add(
    $wsZh.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11347c2a6e2f31d8ba59a034b0f5b586dfe56f2d/e2e/e0757593-76dd-44f0-b232-3cdbe96e3da1.md",
    "",
    "",
    "e0757593-76dd-44f0-b232-3cdbe96e3da1.md"
) | Out-Null
$wsZh.Range("I3").Font.Color = 15570276
$wsZh.Range("J3").Value = "e0757593-76dd-44f0-b232-3cdbe96e3da1.ed61f209bf65f9c275ed609297e9eee89058bdf1.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-06 02:53:23"

Set-RawColumnWidth $wsZh.Columns.Item(3) 29.9777047293527
Set-RawColumnWidth $wsZh.Columns.Item(9) 40
Set-RawColumnWidth $wsZh.Columns.Item(10) 40

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

# Row 2 (7d8b69fb...)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11347c2a6e2f31d8ba59a034b0f5b586dfe56f2d/e2e/7d8b69fb-74ad-419a-b316-fb9822822b42.md",
    "",
    "",
    "7d8b69fb-74ad-419a-b316-fb9822822b42.md"
) | Out-Null
$wsDe.Range("I2").Font.Color = 15570276
$wsDe.Range("J2").Value = "2016-09-06 02:53:30"
$wsDe.Range("K2").Value = "7d8b69fb-74ad-419a-b316-fb9822822b42.e7ae820631db705b393e5b3c6834a1c08189e120.de-de.xlf"

# Row 3 (e0757593...)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11347c2a6e2f31d8ba59a034b0f5b586dfe56f2d/e2e/e0757593-76dd-44f0-b232-3cdbe96e3da1.md",
    "",
    "",
    "e0757593-76dd-44f0-b232-3cdbe96e3da1.md"
) | Out-Null
$wsDe.Range("I3").Font.Color = 15570276
$wsDe.Range("J3").Value = "e0757593-76dd-44f0-b232-3cdbe96e3da1.ed61f209bf65f9c275ed609297e9eee89058bdf1.de-de.xlf"
$wsDe.Range("K3").Value = "7d8b69fb-74ad-419a-b316-fb9822822b42.e7ae820631db705b393e5b3c6834a1c08189e120.de-de.xlf"

Set-RawColumnWidth $wsDe.Columns.Item(3) 29.9777047293527
Set-RawColumnWidth $wsDe.Columns.Item(9) 40
Set-RawColumnWidth $wsDe.Columns.Item(10) 40
